$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Expand the abbreviated correct_ans codes in column L ("b"/"r"/"y")
#     into their full words ("center"/"right"/"left"). ---
$ansMap = @{ "b" = "center"; "r" = "right"; "y" = "left" }

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 12)   # column L = correct_ans
    $val = $cell.Value2
    if ($null -ne $val -and $ansMap.ContainsKey($val)) {
        $cell.Value = $ansMap[$val]
    }
}

# --- 2) Rename the "face" stimulus category to "book" everywhere it
#     appears in the image-path columns (A-D hold promptFile/correctFile/
#     dist_01File/dist_02File entries like "face//face_20.jpg"). ---
$lastCol = $used.Columns.Count

for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le 4; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($null -ne $val -and $val -like "face*") {
            $cell.Value = ($val -replace "face", "book")
        }
    }
}
